# Fixing errors in example upload files.
# - Add a missing practitioner row (P01 @ PHN999:NFP02) to the "Practitioners" sheet
# - Adjust the active selection on "Service Contacts" and "Practitioners"
# - Tidy up a couple of column widths that Excel re-flowed while editing

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Practitioners sheet: append row 6 with data for practitioner P01 under the
# second organisation (PHN999:NFP02), mirroring the existing P01 row already
# present for PHN999:NFP01.
# ---------------------------------------------------------------------------
$wsPrac = $wb.Worksheets.Item("Practitioners")

$wsPrac.Range("A6").Value = "PHN999:NFP02"
$wsPrac.Range("B6").Value = "P01"
$wsPrac.Range("C6").Value = 8
$wsPrac.Range("D6").Value = 1
$wsPrac.Range("E6").Value = 1973
$wsPrac.Range("F6").Value = 2
$wsPrac.Range("G6").Value = 1
$wsPrac.Range("H6").Value = 1
$wsPrac.Range("I6").Value = "tag1"

# Re-fit a few columns (width nudges left over from editing the sheet).
$wsPrac.Columns.Item(1).ColumnWidth = 13.833333333333334
$wsPrac.Columns.Item(3).ColumnWidth = 12.166666666666666
$wsPrac.Columns.Item(6).ColumnWidth = 12

# Move the selection/active cell to column G (whole column selected).
$wsPrac.Activate() | Out-Null
$wsPrac.Range("G1:G1048576").Select() | Out-Null

# ---------------------------------------------------------------------------
# Service Contacts sheet: widen column A and move the selection to D3.
# ---------------------------------------------------------------------------
$wsSvc = $wb.Worksheets.Item("Service Contacts")

$wsSvc.Columns.Item(1).ColumnWidth = 13.666666666666666

$wsSvc.Activate() | Out-Null
$wsSvc.Range("D3").Select() | Out-Null

# ---------------------------------------------------------------------------
# Restore the original active sheet/selection (Metadata!C2) so we don't leave
# the workbook pointed at a different tab than before.
# ---------------------------------------------------------------------------
$wsMeta = $wb.Worksheets.Item("Metadata")
$wsMeta.Activate() | Out-Null
$wsMeta.Range("C2").Select() | Out-Null
